# Insert a new data row at row 41 (pushing existing rows 41..140 down to 42..141)
# and populate it with a new price record dated 2023-11-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(41).Insert()

$ws.Cells.Item(41, 1).Value = 10
$ws.Cells.Item(41, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(41, 3).Value = "La Araucanía"
$ws.Cells.Item(41, 4).Value = (Get-Date -Year 2023 -Month 11 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(41, 5).Value = 9
$ws.Cells.Item(41, 6).Value = 300000001
$ws.Cells.Item(41, 7).Value = "Rabanito"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 45
$ws.Cells.Item(41, 11).Value = 9000
$ws.Cells.Item(41, 12).Value = 9000
$ws.Cells.Item(41, 13).Value = 9000
$ws.Cells.Item(41, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(41, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(41, 16).Value = 750
$ws.Cells.Item(41, 17).Value = 12
$ws.Cells.Item(41, 18).Value = "Hortaliza"
